$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Steps")

# Add the new "Weights" column (F) with weight values for the two polls.
$ws.Range("F1").Value = "Weights"
$ws.Range("F2").Value = "0.9, 1.1"
$ws.Range("F5").Value = "0.9, 1.0, 1.1, 1.2"

# Match the selection left behind in the authored workbook.
$ws.Range("F5").Select()
